$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of this block (rows 60-61), pushing the
# existing rows 60-74 down to 62-76.
$ws.Rows("60:61").Insert()

# New row 60: Sandia / Extra, week of 44559
$ws.Range("A60").Value = 8
$ws.Range("B60").Value = "Terminal La Palmera de La Serena"
$ws.Range("C60").Value = "Coquimbo"
$ws.Range("D60").Value = 44559
$ws.Range("E60").Value = 4
$ws.Range("F60").Value = 100112028
$ws.Range("G60").Value = "Sandia"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Extra"
$ws.Range("J60").Value = 3600
$ws.Range("K60").Value = 3500
$ws.Range("L60").Value = 4000
$ws.Range("M60").Value = 3750
$ws.Range("N60").Value = "$/unidad"
$ws.Range("O60").Value = "Región de O'Higgins"
$ws.Range("P60").Value = 3750
$ws.Range("Q60").Value = 1
$ws.Range("R60").Value = "Hortaliza"

# New row 61: Sandia / Primera, week of 44559
$ws.Range("A61").Value = 8
$ws.Range("B61").Value = "Terminal La Palmera de La Serena"
$ws.Range("C61").Value = "Coquimbo"
$ws.Range("D61").Value = 44559
$ws.Range("E61").Value = 4
$ws.Range("F61").Value = 100112028
$ws.Range("G61").Value = "Sandia"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 2500
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = 2750
$ws.Range("N61").Value = "$/unidad"
$ws.Range("O61").Value = "Región de O'Higgins"
$ws.Range("P61").Value = 2750
$ws.Range("Q61").Value = 1
$ws.Range("R61").Value = "Hortaliza"
